$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabla#1 - Error")

# Insert a new row at row 4, shifting existing rows 4-12 down to 5-13
$ws.Rows.Item(4).Insert()

# Fill in the values for the newly inserted row 4
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = 14
$ws.Cells.Item(4, 4).Value = "@"

# Column A is a static sequential "No" counter, not a formula, so after the
# insert the rows below the new one still carry their old numbers. Renumber
# them to stay sequential (1..12) down through the now-last row 13.
for ($r = 5; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
